$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a single empty "template" data row (row 7), an
# empty totals row (row 8, merged P8:Q8) and a footer row (row 9). We need
# three filled data rows, a totals row with a value, and the footer pushed
# down two rows (to row 11) with a refreshed timestamp.

# Insert two new rows above the current row 8 so the existing row 8
# (totals) becomes row 10 and the existing row 9 (footer) becomes row 11.
$ws.Rows("8:9").Insert()

# Copy the formatting (incl. merged-cell layout) of the template row (7)
# down into the two newly inserted rows (8 and 9).
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q9").PasteSpecial(-4122)

$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

# Columns that hold the textual numbers coming from the POS export must be
# stored as text (not auto-converted to numbers), so force the Text format
# before writing the values - matches the numFmtId 49 change.
$ws.Range("C7:G9").NumberFormat = "@"
$ws.Range("N7:O9").NumberFormat = "@"
$ws.Range("H7:K9").NumberFormat = "@"
$ws.Range("Q7:Q9").NumberFormat = "@"
$ws.Range("L7:M9").NumberFormat = "@"

# Row 7 - item 1
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "BOBAI SUNSCREEN TINTED CREAM SPF 80 - 60 GM"
$ws.Range("H7").Value = "0:0"
$ws.Range("L7").Value = "1"
$ws.Range("N7").Value = "195.00"
$ws.Range("P7").Value = "195.0000"
$ws.Range("Q7").Value = "1:0"

# Row 8 - item 2
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "KADEE 4000 CAPS"
$ws.Range("H8").Value = "0:0"
$ws.Range("L8").Value = "0"
$ws.Range("N8").Value = "180.00"
$ws.Range("P8").Value = "180.0000"
$ws.Range("Q8").Value = "1:0"

# Row 9 - item 3
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "VONSECA 20MG 14 TAB"
$ws.Range("H9").Value = "0:1"
$ws.Range("L9").Value = "1"
$ws.Range("N9").Value = "140.00"
$ws.Range("P9").Value = "140.0000"
$ws.Range("Q9").Value = "1:0"

# Row 10 - totals (number of transactions)
$ws.Range("P10").Value = 515

# Row 11 - footer, refreshed timestamp
$ws.Range("A11").Value = "Tuesday, 30 September, 2025 9:27 AM"
